$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update sort/lfpExported/spkExported flag cells from "0" to "1" ---
# Done via copy/paste-special (value+formats) from a same-styled source cell
# that already holds the shared string "1", so the written cell keeps the
# exact original style and becomes a text "1" (not a numeric 1).
$ws.Range("H3").Copy()
$ws.Range("I78").PasteSpecial(-4104)
$ws.Range("J78").PasteSpecial(-4104)
$ws.Range("I79").PasteSpecial(-4104)
$ws.Range("J79").PasteSpecial(-4104)
$ws.Range("I80").PasteSpecial(-4104)
$ws.Range("J80").PasteSpecial(-4104)
$ws.Range("I81").PasteSpecial(-4104)
$ws.Range("J81").PasteSpecial(-4104)
$ws.Range("I82").PasteSpecial(-4104)
$ws.Range("J82").PasteSpecial(-4104)
$ws.Range("H118").PasteSpecial(-4104)
$ws.Range("H119").PasteSpecial(-4104)
$ws.Range("H120").PasteSpecial(-4104)
$ws.Range("H122").PasteSpecial(-4104)
$ws.Range("H123").PasteSpecial(-4104)
$ws.Range("H124").PasteSpecial(-4104)
$ws.Range("H125").PasteSpecial(-4104)
$ws.Range("H126").PasteSpecial(-4104)
$ws.Range("H127").PasteSpecial(-4104)
$ws.Range("H129").PasteSpecial(-4104)
$ws.Range("H130").PasteSpecial(-4104)
$ws.Range("H131").PasteSpecial(-4104)
$ws.Range("H133").PasteSpecial(-4104)
$ws.Range("H134").PasteSpecial(-4104)
$ws.Range("H135").PasteSpecial(-4104)
$ws.Range("H137").PasteSpecial(-4104)
$ws.Range("H138").PasteSpecial(-4104)
$ws.Range("H139").PasteSpecial(-4104)
$ws.Range("H140").PasteSpecial(-4104)
$ws.Range("H141").PasteSpecial(-4104)
$ws.Range("H143").PasteSpecial(-4104)
$ws.Range("H144").PasteSpecial(-4104)
$ws.Range("H145").PasteSpecial(-4104)
$ws.Range("H146").PasteSpecial(-4104)
$ws.Range("H147").PasteSpecial(-4104)
$ws.Range("H149").PasteSpecial(-4104)
$ws.Range("H150").PasteSpecial(-4104)
$ws.Range("H151").PasteSpecial(-4104)
$ws.Range("H152").PasteSpecial(-4104)
$ws.Range("H153").PasteSpecial(-4104)
$ws.Range("H155").PasteSpecial(-4104)
$ws.Range("H156").PasteSpecial(-4104)
$ws.Range("H157").PasteSpecial(-4104)
$ws.Range("H158").PasteSpecial(-4104)
$ws.Range("H159").PasteSpecial(-4104)
$ws.Range("H161").PasteSpecial(-4104)
$ws.Range("H162").PasteSpecial(-4104)
$ws.Range("H163").PasteSpecial(-4104)
$ws.Range("H164").PasteSpecial(-4104)
$ws.Range("H165").PasteSpecial(-4104)
$ws.Range("H167").PasteSpecial(-4104)
$ws.Range("H168").PasteSpecial(-4104)
$ws.Range("H169").PasteSpecial(-4104)
$ws.Range("H170").PasteSpecial(-4104)
$ws.Range("H171").PasteSpecial(-4104)
$ws.Range("H173").PasteSpecial(-4104)
$ws.Range("H174").PasteSpecial(-4104)
$ws.Range("H175").PasteSpecial(-4104)
$ws.Range("H176").PasteSpecial(-4104)

$ws.Range("H22").Copy()
$ws.Range("H121").PasteSpecial(-4104)
$ws.Range("H128").PasteSpecial(-4104)
$ws.Range("H136").PasteSpecial(-4104)
$ws.Range("H142").PasteSpecial(-4104)
$ws.Range("H148").PasteSpecial(-4104)
$ws.Range("H154").PasteSpecial(-4104)
$ws.Range("H160").PasteSpecial(-4104)
$ws.Range("H166").PasteSpecial(-4104)

$ws.Range("H47").Copy()
$ws.Range("H132").PasteSpecial(-4104)
$ws.Range("H172").PasteSpecial(-4104)
$ws.Range("H177").PasteSpecial(-4104)

# --- Update the view / selection state on Sheet1 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 112
$win.ScrollColumn = 1
$ws.Range("D138").Select()

# --- Update the workbook window placement/size ---
$win.Left = -120
$win.Top = -120
$win.Width = 29040
$win.Height = 15840
